$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row (row 20) below the last existing worker row (19),
#    shifting everything below (the signature block) down by one row.
# ---------------------------------------------------------------------------
$ws.Range("B20:J20").Insert(-4121)   # xlShiftDown

# ---------------------------------------------------------------------------
# 2) The old last row (19) had the "closing" border style (bottom border of
#    the table). After inserting a new row below it, row 19 becomes a
#    regular interior row (same look as row 18) and the new row 20 becomes
#    the new closing row (same look the old row 19 had).
# ---------------------------------------------------------------------------

# Give the new row 20 the borders/format that row 19 used to have (closing row)
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Give row 19 the borders/format of an interior row (same as row 18)
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the values for the new row (same worker, new period 2509)
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = $ws.Range("B19").Value2
$ws.Range("C20").Value = $ws.Range("C19").Value2
$ws.Range("D20").Value = $ws.Range("D19").Value2
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# ---------------------------------------------------------------------------
# 4) Update the summary figures: one more overdue period (4 -> 5) and the
#    new total overdue amount (227760 + 56940 = 284700)
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = 5
$ws.Range("E11").Value = 284700

Write-Host "Edit complete"
